$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: pull the paragraph's own opening "<w:p ...>" tag (attributes like
# w14:paraId / w:rsidR / w:rsidP) out of its WordOpenXML so a replacement can
# keep them instead of minting a bare "<w:p>".
# ---------------------------------------------------------------------------
function Get-ParaOpenTag($para) {
    $xml = $para.Range.WordOpenXML
    if ($xml -match '(<w:p[ >][^>]*>|<w:p/>)') {
        $tag = $matches[0]
        if ($tag.EndsWith('/>')) {
            $attrs = $tag.Substring(4, $tag.Length - 6)
            return "<w:p" + $attrs + ">"
        } else {
            return $tag
        }
    }
    return "<w:p>"
}

function New-PackageXml($bodyInnerXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
        '<w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Replace-ParagraphRuns($para, $runsXml) {
    $openTag = Get-ParaOpenTag($para)
    $r = $para.Range
    $r.End = $r.End - 1
    if ($r.Start -lt $r.End) {
        $r.Delete()
    }
    $insertPoint = $para.Range
    $insertPoint.End = $insertPoint.End - 1
    $bodyXml = $openTag + $runsXml + "</w:p>"
    $insertPoint.InsertXML((New-PackageXml $bodyXml))
}

# ---------------------------------------------------------------------------
# Locate the two target paragraphs by their (pre-edit) text.
# ---------------------------------------------------------------------------
$target1 = $Null
$linked = $Null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt -eq "Secondary Storage (File I/O)`r") {
        $target1 = $para
    }
    if ($txt -eq "Linked Lists`r") {
        $linked = $para
    }
}

# ---------------------------------------------------------------------------
# Change 1 (Week 8 section): merge "Secondary Storage" / " (File I/O)" into
# one run, then add a new run " - Text Files".
# ---------------------------------------------------------------------------
$runs1 = '<w:r><w:t>Secondary Storage (File I/O)</w:t></w:r><w:r><w:t xml:space="preserve"> &#8211; Text Files</w:t></w:r>'
Replace-ParagraphRuns $target1 $runs1

# ---------------------------------------------------------------------------
# Change 2 (Week 12 section): the empty paragraph right after "Linked Lists"
# becomes "Secondary Storage (File I/O) - Records and Fields" (two runs).
# ---------------------------------------------------------------------------
$target2 = $linked.Next()
$runs2 = '<w:r><w:t xml:space="preserve">Secondary Storage (File I/O) &#8211; </w:t></w:r><w:r><w:t>Records and Fields</w:t></w:r>'
Replace-ParagraphRuns $target2 $runs2

Write-Output "done"
